$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.316.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.28%  "

$ws.Range("D3").Value = "'1.792.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.16%  "

$ws.Range("E4").Value = "  +0.77%  "

$ws.Range("D5").Value = "'1.009"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.69%  "

$ws.Range("D6").Value = "'307.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.07%  "

$ws.Range("D7").Value = "'0.4533"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.11%  "

$ws.Range("D8").Value = "'0.3606"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.45%  "

$ws.Range("D9").Value = "'0.07067"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.67%  "

$ws.Range("D10").Value = "'0.8715"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.78%  "

$ws.Range("E11").Value = "  -0.87%  "

$ws.Range("D12").Value = "'19.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.44%  "

$ws.Range("D13").Value = "'1.836.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.65%  "

$ws.Range("D14").Value = "'5.257"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.42%  "

$ws.Range("D15").Value = "'6.298"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.13%  "

$ws.Range("D16").Value = "'84.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.71%  "

$ws.Range("D17").Value = "'1.011"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.75%  "

$ws.Range("D18").Value = "'0.000008522"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.74%  "

$ws.Range("D19").Value = "'1.009"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.67%  "

$ws.Range("D20").Value = "'26.409.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.08%  "

$ws.Range("D21").Value = "'14.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.98%  "

$ws.Range("D22").Value = "'4.968"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.84%  "

$ws.Range("D23").Value = "'2.066.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.04%  "

$ws.Range("D24").Value = "'10.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.81%  "

$ws.Range("D25").Value = "'1.984"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.42%  "

$ws.Range("D26").Value = "'151.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.17%  "

$ws.Range("D27").Value = "'17.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.59%  "

$ws.Range("D28").Value = "'2.035"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.54%  "

$ws.Range("D29").Value = "'112.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.25%  "

$ws.Range("D30").Value = "'4.814"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.87%  "

$ws.Range("D31").Value = "'0.08654"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.22%  "

$ws.Range("D32").Value = "'3.028"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.20%  "

$ws.Range("D33").Value = "'4.440"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.62%  "

$ws.Range("D34").Value = "'0.7144"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.26%  "

$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").Value = "'2.638"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.84%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.095"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.63%  "

$ws.Range("D37").Value = "'1.009"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.58%  "

$ws.Range("D38").Value = "'1.080"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.39%  "

$ws.Range("D39").Value = "'0.01944"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("D40").Value = "'0.05086"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.36%  "

$ws.Range("D41").Value = "'2.875"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.73%  "

$ws.Range("D42").Value = "'6.903"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.46%  "

$ws.Range("D43").Value = "'0.4935"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.89%  "

$ws.Range("D44").Value = "'0.1524"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.63%  "

$ws.Range("D45").Value = "'7.942"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.92%  "

$ws.Range("E46").Value = "  +0.76%  "

$ws.Range("D47").Value = "'0.4554"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.12%  "

$ws.Range("D48").Value = "'9.868"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.54%  "

$ws.Range("D49").Value = "'100.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.73%  "

$ws.Range("D50").Value = "'1.575"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.36%  "

$ws.Range("D51").Value = "'0.05951"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.03%  "
